$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per-row and are being permuted: D, M, N, O, P, Q, R, S, T
$cols = @("D","M","N","O","P","Q","R","S","T")

# Snapshot the original values for each affected row/column before any writes,
# since the change is a permutation of rows (some destination rows pull data
# from rows that would otherwise get overwritten earlier in the loop).
# Value2 is used (rather than Value) to get/set plain scalars reliably.
$original = @{}
for ($r = 2; $r -le 26; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Map: destination row -> source row (the destination row receives the ORIGINAL
# data that the source row held before this script ran)
$rowMap = @{
    2 = 20
    3 = 8
    4 = 16
    5 = 12
    6 = 4
    7 = 21
    8 = 15
    9 = 13
    10 = 25
    11 = 22
    12 = 18
    13 = 11
    14 = 17
    15 = 23
    16 = 14
    17 = 26
    18 = 19
    19 = 5
    20 = 6
    21 = 24
    22 = 7
    23 = 2
    24 = 10
    25 = 3
    26 = 9
}

foreach ($dst in $rowMap.Keys) {
    $src = $rowMap[$dst]
    $srcVals = $original[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dst").Value2 = $srcVals[$c]
    }
}
